$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A17").Value = "spese annuali"
$ws.Range("B17").Value = "non salva nel db le note e non midifica colore del pannello"
$ws.Range("C17").Value = "verificare inserimento"
$ws.Range("D17").Value = "da verificare"
$ws.Range("E17").Value = 44883
$ws.Range("E17").NumberFormat = $ws.Range("E16").NumberFormat

$ws.Range("A18").Select()
